$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Column A width (XML width=23 -> ColumnWidth is offset by ~0.8333 in this engine)
$ws.Columns.Item(1).ColumnWidth = 22.166666666666668

# Row 2 - maize silage
$ws.Range("B2").Value = 40.085336
$ws.Range("C2").Value = 1.540789
$ws.Range("D2").Value = 0.817209

# Row 3 - grass silage
$ws.Range("B3").Value = 36.613181
$ws.Range("C3").Value = 2.473993
$ws.Range("D3").Value = 0.780554

# Row 4 - sugar beet silage
$ws.Range("B4").Value = 49.423171
$ws.Range("C4").Value = 0.560114
$ws.Range("D4").Value = 0.06216

# Row 5 - cattle manure
$ws.Range("B5").Value = 5.378375
$ws.Range("C5").Value = 0.778969
$ws.Range("D5").Value = 0.205133

# Row 6 - remove the "cattle manure very uncertain" row data (keep row, clear content)
$ws.Range("A6:D6").ClearContents()

# Apply number format "0.000" to B2:D5
$ws.Range("B2:D5").NumberFormat = "0.000"

# Selection
$ws.Range("C7").Select()
